$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data scraped by GitHub Actions
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.306.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.494.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.493.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.084.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.489.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.328.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -9.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.565"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.631.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.52%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.513.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "23.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "161.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0781"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.473.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.50%  "
